$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Starting layout (rows 1-16):
#   Row 1             : header  (PreCondName | TrgCol)
#   Rows 2-6  (5)     : CS+3.BMP / -0.1
#   Rows 7-11 (5)     : CS+1.BMP / 0.4
#   Rows 12-16(5)     : CS+4.BMP / 0.4
#
# Target layout (rows 1-25):
#   Row 1             : header, unchanged
#   Rows 2-9  (8)     : CS+3.BMP / -0.1   -- rows 6-8 use a NEW highlighted
#                        (light-grey) fill style; row 9 reverts to the
#                        original style
#   Rows 10-17(8)     : CS+1.BMP / 0.4
#   Rows 18-25(8)     : CS+4.BMP / 0.4
# ---------------------------------------------------------------------------

# 1) Grow each 5-row block to 8 rows by inserting 3 fresh rows right after
#    each block. Excel's row Insert naturally copies formatting down from
#    the row immediately above, which is exactly what we want here.

# CS+3 block: rows 2-6 -> insert after row 6 (new rows become 7,8,9)
$ws.Rows("7:9").Insert()

# CS+1 block: was rows 7-11, now shifted to rows 10-14 -> insert after
# (new rows become 15,16,17)
$ws.Rows("15:17").Insert()

# CS+4 block: was rows 12-16, now shifted twice (+3, +3) to rows 18-22 ->
# insert after (new rows become 23,24,25)
$ws.Rows("23:25").Insert()

# 2) Fill in the values for every newly inserted (currently blank) cell.

# CS+3 rows 7,8,9
for ($r = 7; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = "Stimuli/CS+3.BMP"
    $ws.Cells.Item($r, 2).Value = -0.1
}

# CS+1 rows 15,16,17
for ($r = 15; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = "Stimuli/CS+1.BMP"
    $ws.Cells.Item($r, 2).Value = 0.4
}

# CS+4 rows 23,24,25
for ($r = 23; $r -le 25; $r++) {
    $ws.Cells.Item($r, 1).Value = "Stimuli/CS+4.BMP"
    $ws.Cells.Item($r, 2).Value = 0.4
}

# 3) Rows 6,7,8 (A & B) get a new light-grey highlight fill, distinguishing
#    the last 3 "Pre-Cond" trials from the rest of the CS+3 block. Build
#    the style once on row 6, then replicate it onto rows 7 and 8.
$hl = $ws.Range("A6:B6")
$hl.Interior.Color = 15921906   # RGB(242,242,242) -> 0xF2F2F2

$ws.Range("A6:B6").Copy($ws.Range("A7:B7"))
$ws.Range("A6:B6").Copy($ws.Range("A8:B8"))

# Row 9 explicitly keeps the ORIGINAL (non-highlighted) CS+3 style, matching
# rows 2-5. Re-copy formatting from row 5 onto row 9 just to be certain it
# wasn't touched by anything above (defensive / explicit per the diff).
$ws.Range("A5:B5").Copy($ws.Range("A9:B9"))
$ws.Cells.Item(9, 1).Value = "Stimuli/CS+3.BMP"
$ws.Cells.Item(9, 2).Value = -0.1

# 4) Move the active selection to A31, matching the saved selection.
$ws.Range("A31").Select()

Write-Output "done"
Write-Output $ws.UsedRange.Address()
